$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.494.00"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "3.432.47"
$ws.Range("E3").Value = "  -5.28%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'570.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.11%  "
$ws.Range("D6").Value = "'188.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("E7").Value = "  -3.01%  "
$ws.Range("D8").Value = "3.425.44"
$ws.Range("E8").Value = "  -5.14%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").Value = "'0.609"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("D12").Value = "'50.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  -7.47%  "
$ws.Range("D14").Value = "'8.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.90%  "
$ws.Range("D15").Value = "3.979.25"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("D16").Value = "'627.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.66%  "
$ws.Range("D17").Value = "68.323.02"
$ws.Range("E17").Value = "  -3.04%  "
$ws.Range("D18").Value = "3.445.31"
$ws.Range("E18").Value = "  -4.98%  "
$ws.Range("D20").Value = "'12.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.44%  "
$ws.Range("D21").Value = "'17.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.60%  "
$ws.Range("D22").Value = "'0.932"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.66%  "
$ws.Range("D23").Value = "'17.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").Value = "'5.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "'98.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.11%  "
$ws.Range("E26").Value = "  -8.46%  "
$ws.Range("E27").Value = "  -5.63%  "
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "'9.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.84%  "
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("D31").Value = "'32.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D32").Value = "'4.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.72%  "
$ws.Range("E33").Value = "  -8.57%  "
$ws.Range("D34").Value = "'11.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.61%  "
$ws.Range("E35").Value = "  -7.46%  "
$ws.Range("D36").Value = "'60.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "3.618.48"
$ws.Range("E38").Value = "  -8.31%  "
$ws.Range("D39").Value = "0.0₃0768"
$ws.Range("E39").Value = "  -13.29%  "
$ws.Range("D40").Value = "'494.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.50%  "
$ws.Range("D41").Value = "'3.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  -6.82%  "
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D43").Value = "'3.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +74.22%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.364"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.19%  "
$ws.Range("D45").Value = "'0.131"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").Value = "'33.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.47%  "
$ws.Range("D47").Value = "'0.0432"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.07%  "
$ws.Range("D48").Value = "'3.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.73%  "
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("E50").Value = "  -5.17%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
